$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-10 Friday" "2025-10-11 Saturday"

Replace-Text "886÷8=110, 6" "420÷5=84, 0"
Replace-Text "952÷5=190, 2" "190÷6=31, 4"
Replace-Text "462÷4=115, 2" "320÷9=35, 5"
Replace-Text "140÷4=35, 0" "366÷2=183, 0"
Replace-Text "578÷5=115, 3" "577÷5=115, 2"

Replace-Text "168÷5=33, 3" "450÷4=112, 2"
Replace-Text "964÷7=137, 5" "762÷7=108, 6"
Replace-Text "406÷4=101, 2" "946÷7=135, 1"
Replace-Text "194÷8=24, 2" "387÷4=96, 3"
Replace-Text "185÷4=46, 1" "442÷4=110, 2"

Replace-Text "882÷9=98, 0" "388÷3=129, 1"
Replace-Text "891÷5=178, 1" "723÷5=144, 3"
Replace-Text "108÷7=15, 3" "584÷7=83, 3"
Replace-Text "484÷6=80, 4" "643÷4=160, 3"
Replace-Text "468÷2=234, 0" "183÷9=20, 3"

Replace-Text "899÷8=112, 3" "133÷3=44, 1"
Replace-Text "593÷9=65, 8" "781÷2=390, 1"
Replace-Text "170÷5=34, 0" "981÷6=163, 3"
Replace-Text "615÷3=205, 0" "154÷3=51, 1"
Replace-Text "771÷6=128, 3" "212÷3=70, 2"

Replace-Text "931÷4=232, 3" "525÷9=58, 3"
Replace-Text "650÷4=162, 2" "891÷2=445, 1"
Replace-Text "707÷3=235, 2" "671÷6=111, 5"
Replace-Text "257÷9=28, 5" "118÷2=59, 0"
Replace-Text "498÷8=62, 2" "877÷2=438, 1"

Write-Output "Done"
